# Adds Barnsley's organisational code to the project_names lookup table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("project_names")

# The table (Table1) currently spans A1:G31 with data rows 2-31.
# Add a new row 32 for Barnsley, mirroring the existing Rotherham-style row
# that shares the "Barnsley Bassetlaw Rotherham" project_name_place grouping.
$newRow = 32

$ws.Range("A$newRow").Value = "Barnsley"
$ws.Range("B$newRow").Value = "02P00"
$ws.Range("C$newRow").Value = "Barnsley Bassetlaw Rotherham"
$ws.Range("D$newRow").Value = "0"

# Expand the table (ListObject) to include the newly added row.
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:G$newRow"))

$ws.Range("C$newRow").Select()

$wb.Save()
